$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "Main" sheet to "About"
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("Main")
$wsAbout.Name = "About"

# ---------------------------------------------------------------------------
# 2) "About" sheet: insert 8 explanatory rows right after the "Notes:" header
#    (old row 10 becomes new row 18), fill in the new leakage-rate
#    explanation text, and fix up the hyperlink that shifted down.
# ---------------------------------------------------------------------------
$wsAbout.Rows("10:17").Insert()

$wsAbout.Range("A10").Value = "For every unit of CO2e not emitted in this country as a result of the modeled policy package, a fraction"
$wsAbout.Range("A11").Value = "of a unit of CO2e is emitted in a foreign country. This may be due to businesses relocating or opening"
$wsAbout.Range("A12").Value = "elsewhere rather than in this country, or it may be due to changes in production levels of existing"
$wsAbout.Range("A13").Value = "businesses. These leakage rates were determined via a computer model based on a carbon tax, so"
$wsAbout.Range("A14").Value = "leakage can be negative if a foreign firm has a lower carbon intensity than a U.S. firm when"
$wsAbout.Range("A15").Value = "producing a particular good (this is the case for natural gas and coal), or if imports to the U.S. are"
$wsAbout.Range("A16").Value = "reduced (this is the case for oil)."

$wsAbout.Range("A10:A17").Font.Bold = $false

# Fix the hyperlinks: row insertion does not auto-shift the Hyperlinks
# collection, so rebuild both entries at their correct (shifted) locations,
# then restore the Hyperlink cell style that Hyperlinks.Add() perturbs.
$wsAbout.Hyperlinks.Delete() | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("B6"), "http://www.rff.org/RFF/Documents/RFF-DP-10-47.pdf", "") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("A22"), "http://www.worldcement.com/news/cement/articles/Cement_global_trading_patterns_961.aspx", ".UyvL5fldWE8") | Out-Null
$wsAbout.Range("B6").Style = "Hyperlink"
$wsAbout.Range("A22").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3) "FLRbI" sheet: clarify the "Leakage Rate" header to note it is
#    dimensionless.
# ---------------------------------------------------------------------------
$wsFLRbI = $wb.Worksheets.Item("FLRbI")
$wsFLRbI.Range("B1").Value = "Leakage Rate (dimensionless)"
$wsFLRbI.Range("B2").Select() | Out-Null

# Leave "About" as the active sheet/selection, matching the saved workbook.
$wsAbout.Activate() | Out-Null
$wsAbout.Range("B33").Select() | Out-Null
